# Fixed rota print date bug
# Update the weekly rota names in column B (and the new helper entries in
# column C) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "MahaDeva PM"
$ws.Range("B3").Value = "MahaDeva - OHC"
$ws.Range("C3").Value = "MahaDeva - OHC"

$ws.Range("C4").Value = "Adam"

$ws.Range("B6").Value = "MahaDeva AM"

$ws.Range("B10").Value = "Dganit"
$ws.Range("B11").Value = "Shakti"
$ws.Range("B12").Value = "Anuka"
$ws.Range("B13").Value = "Ben"
$ws.Range("B14").Value = "Mahi"
$ws.Range("B15").Value = "RAP"
